$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.219.39"
$ws.Range("E2").Value = "  -1.02%  "

$ws.Range("D3").Value = "3.531.32"
$ws.Range("E3").Value = "  +0.50%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "607.27"
$ws.Range("E5").Value = "  -0.20%  "

$ws.Range("D6").Value = "143.34"
$ws.Range("E6").Value = "  -3.16%  "

$ws.Range("D7").Value = "3.529.14"
$ws.Range("E7").Value = "  +0.46%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "0.480"
$ws.Range("E9").Value = "  +0.44%  "

$ws.Range("E10").Value = "  -4.27%  "

$ws.Range("E11").Value = "  +0.74%  "

$ws.Range("E12").Value = "  -2.90%  "

$ws.Range("D13").Value = "4.131.15"
$ws.Range("E13").Value = "  +0.55%  "

$ws.Range("D14").Value = "0.0000208"
$ws.Range("E14").Value = "  -4.60%  "

$ws.Range("D15").Value = "30.14"
$ws.Range("E15").Value = "  -5.53%  "

$ws.Range("D16").Value = "3.528.32"
$ws.Range("E16").Value = "  +0.49%  "

$ws.Range("D17").Value = "66.311.30"
$ws.Range("E17").Value = "  -1.02%  "

$ws.Range("E18").Value = "  -0.69%  "

$ws.Range("D19").Value = "10.93"
$ws.Range("E19").Value = "  +2.07%  "

$ws.Range("E20").Value = "  -4.02%  "

$ws.Range("D21").Value = "14.93"
$ws.Range("E21").Value = "  -2.70%  "

$ws.Range("D22").Value = "425.26"
$ws.Range("E22").Value = "  -2.98%  "

$ws.Range("D23").Value = "0.601"
$ws.Range("E23").Value = "  -1.34%  "

$ws.Range("D24").Value = "78.63"
$ws.Range("E24").Value = "  -0.98%  "

$ws.Range("D25").Value = "3.674.39"
$ws.Range("E25").Value = "  +0.50%  "

$ws.Range("E27").Value = "  -1.39%  "

$ws.Range("E28").Value = "  -2.71%  "

$ws.Range("D29").Value = "9.16"
$ws.Range("E29").Value = "  -6.22%  "

$ws.Range("D30").Value = "2.48"
$ws.Range("E30").Value = "  -1.72%  "

$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("E32").Value = "  -4.62%  "

$ws.Range("E33").Value = "  -6.34%  "

$ws.Range("D34").Value = "25.26"
$ws.Range("E34").Value = "  -1.07%  "

$ws.Range("D35").Value = "3.522.18"
$ws.Range("E35").Value = "  +0.43%  "

$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("D37").Value = "1.76"
$ws.Range("E37").Value = "  -2.96%  "

$ws.Range("D38").Value = "7.80"
$ws.Range("E38").Value = "  -2.94%  "

$ws.Range("E39").Value = "  -5.65%  "

$ws.Range("E40").Value = "  -0.06%  "

$ws.Range("D41").Value = "171.78"
$ws.Range("E41").Value = "  -0.90%  "

$ws.Range("E42").Value = "  -4.16%  "

$ws.Range("E43").Value = "  -4.31%  "

$ws.Range("E44").Value = "  -0.46%  "

$ws.Range("E45").Value = "  -8.53%  "

$ws.Range("D46").Value = "45.41"
$ws.Range("E46").Value = "  -1.62%  "

$ws.Range("D47").Value = "26.05"
$ws.Range("E47").Value = "  -6.26%  "

$ws.Range("E48").Value = "  -4.86%  "

$ws.Range("E49").Value = "  -2.46%  "

$ws.Range("D50").Value = "7.14"
$ws.Range("E50").Value = "  -4.37%  "

$ws.Range("D51").Value = "0.945"
$ws.Range("E51").Value = "  -4.71%  "
